$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Part number for the 5x2 pin header changed (superseded part) - update BOM
$ws.Range("C13").Style = "Normal"
$ws.Range("C13").Value = "649-68602-110HLF"

# Leave selection where the author left it when saving
$null = $ws.Range("C7").Select()
